# Auto-generated Excel COM-interop script replicating the scheduled runner refresh
# of market/leve profitability figures across the ALC/ARM/BSM/CRP/CUL/GSM/LTW sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 158.5
$ws.Range("I39").Value = 44.666668
$ws.Range("K39").Value = 134.000004
$ws.Range("M39").Value = 161.999996

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1265.6177
$ws.Range("I132").Value = 1257.8572
$ws.Range("K132").Value = 3773.5716
$ws.Range("M132").Value = -1243.5716

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 935.15625
$ws.Range("I135").Value = 884.6774
$ws.Range("J135").Value = 2500
$ws.Range("K135").Value = 7962.096600000001
$ws.Range("L135").Value = 22500
$ws.Range("M135").Value = -5427.096600000001
$ws.Range("N135").Value = -27570

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2424.3333
$ws.Range("I137").Value = 2329.8857
$ws.Range("J137").Value = 2754.9
$ws.Range("K137").Value = 6989.657099999999
$ws.Range("L137").Value = 8264.700000000001
$ws.Range("M137").Value = -4439.657099999999
$ws.Range("N137").Value = -13364.7

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4059.27
$ws.Range("I138").Value = 1938.579
$ws.Range("J138").Value = 4556.716
$ws.Range("K138").Value = 5815.737
$ws.Range("L138").Value = 13670.148
$ws.Range("M138").Value = -675.7370000000001
$ws.Range("N138").Value = -23950.148

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10859.046
$ws.Range("I32").Value = 4526.589
$ws.Range("K32").Value = 4526.589
$ws.Range("M32").Value = -4239.589

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 203148.8
$ws.Range("I61").Value = 2312.6562
$ws.Range("K61").Value = 2312.6562
$ws.Range("M61").Value = -2100.6562

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 5051.675
$ws.Range("I74").Value = 2141.0557
$ws.Range("K74").Value = 2141.0557
$ws.Range("M74").Value = -1267.0557

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 5051.675
$ws.Range("I77").Value = 2141.0557
$ws.Range("K77").Value = 10705.2785
$ws.Range("M77").Value = -6337.2785

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2277.0667
$ws.Range("I132").Value = 1959.225
$ws.Range("J132").Value = 4819.8
$ws.Range("K132").Value = 5877.674999999999
$ws.Range("L132").Value = 14459.4
$ws.Range("M132").Value = -3347.674999999999
$ws.Range("N132").Value = -19519.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 203148.8
$ws.Range("I136").Value = 2312.6562
$ws.Range("K136").Value = 6937.9686
$ws.Range("M136").Value = -4387.9686

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 7162.727
$ws.Range("I86").Value = 10602
$ws.Range("J86").Value = 3035.6
$ws.Range("K86").Value = 10602
$ws.Range("L86").Value = 3035.6
$ws.Range("M86").Value = -9479
$ws.Range("N86").Value = -5281.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 7162.727
$ws.Range("I89").Value = 10602
$ws.Range("J89").Value = 3035.6
$ws.Range("K89").Value = 53010
$ws.Range("L89").Value = 15178
$ws.Range("M89").Value = -47394
$ws.Range("N89").Value = -26410

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 49725.145
$ws.Range("I31").Value = 57111
$ws.Range("J31").Value = 5410
$ws.Range("K31").Value = 57111
$ws.Range("L31").Value = 5410
$ws.Range("M31").Value = -56816
$ws.Range("N31").Value = -6000

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 49725.145
$ws.Range("I34").Value = 57111
$ws.Range("J34").Value = 5410
$ws.Range("K34").Value = 57111
$ws.Range("L34").Value = 5410
$ws.Range("M34").Value = -56909
$ws.Range("N34").Value = -5814

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2761.1333
$ws.Range("I58").Value = 2323.5
$ws.Range("J58").Value = 8888
$ws.Range("K58").Value = 2323.5
$ws.Range("L58").Value = 8888
$ws.Range("M58").Value = -2120.5
$ws.Range("N58").Value = -9294

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 823.7308
$ws.Range("I107").Value = 539.7
$ws.Range("K107").Value = 539.7
$ws.Range("M107").Value = 1380.3

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 12258.917
$ws.Range("I134").Value = 6710.7
$ws.Range("K134").Value = 20132.1
$ws.Range("M134").Value = -17597.1

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2761.1333
$ws.Range("I136").Value = 2323.5
$ws.Range("J136").Value = 8888
$ws.Range("K136").Value = 6970.5
$ws.Range("L136").Value = 26664
$ws.Range("M136").Value = -4420.5
$ws.Range("N136").Value = -31764

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 323.2857
$ws.Range("I8").Value = 323.2857
$ws.Range("K8").Value = 969.8571000000001
$ws.Range("M8").Value = -830.8571000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 9250
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 9250
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 27750
$ws.Range("N80").Value = -29622
$ws.Range("M80").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 9250
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 9250
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 83250
$ws.Range("N83").Value = -92610
$ws.Range("M83").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1285.25
$ws.Range("J113").Value = 1268.375
$ws.Range("L113").Value = 3805.125
$ws.Range("N113").Value = -8145.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 2475.963
$ws.Range("I122").Value = 2474.5386
$ws.Range("J122").Value = 2477.2856
$ws.Range("K122").Value = 22270.8474
$ws.Range("L122").Value = 22295.5704
$ws.Range("M122").Value = -19820.8474
$ws.Range("N122").Value = -27195.5704

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 306.63635
$ws.Range("I2").Value = 235.66667
$ws.Range("K2").Value = 235.66667
$ws.Range("M2").Value = -122.66667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 21163.455
$ws.Range("I70").Value = 19150
$ws.Range("K70").Value = 19150
$ws.Range("M70").Value = -18880

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 21163.455
$ws.Range("I73").Value = 19150
$ws.Range("K73").Value = 19150
$ws.Range("M73").Value = -18214

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1040.2333
$ws.Range("I97").Value = 944.96295
$ws.Range("J97").Value = 1897.6666
$ws.Range("K97").Value = 944.96295
$ws.Range("L97").Value = 1897.6666
$ws.Range("M97").Value = -448.96295
$ws.Range("N97").Value = -2889.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4126.9565
$ws.Range("I7").Value = 3700.7856
$ws.Range("K7").Value = 3700.7856
$ws.Range("M7").Value = -3588.7856

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2625.9033
$ws.Range("I22").Value = 1800.5
$ws.Range("J22").Value = 3305.647
$ws.Range("K22").Value = 1800.5
$ws.Range("L22").Value = 3305.647
$ws.Range("M22").Value = -1505.5
$ws.Range("N22").Value = -3895.647

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 2625.9033
$ws.Range("I27").Value = 1800.5
$ws.Range("J27").Value = 3305.647
$ws.Range("K27").Value = 1800.5
$ws.Range("L27").Value = 3305.647
$ws.Range("M27").Value = -1693.5
$ws.Range("N27").Value = -3519.647

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4138.026
$ws.Range("I40").Value = 3479.5667
$ws.Range("K40").Value = 3479.5667
$ws.Range("M40").Value = -3343.5667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3848.762
$ws.Range("J46").Value = 4114
$ws.Range("L46").Value = 4114
$ws.Range("N46").Value = -4490

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 4126.9565
$ws.Range("I126").Value = 3700.7856
$ws.Range("K126").Value = 11102.3568
$ws.Range("M126").Value = -8632.356800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4644.6
$ws.Range("I132").Value = 3999.8572
$ws.Range("K132").Value = 11999.5716
$ws.Range("M132").Value = -9469.571599999999
